# Trade #83 closed at 2026-02-17 15:53:07 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades, and MarketMaking sheets
# to reflect the newly closed trade #83.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.99   # Current Capital
$summary.Range("B4").Value = -0.02     # Total P&L $
$summary.Range("B5").Value = -0        # Total P&L %
$summary.Range("B6").Value = 83        # Total Trades
$summary.Range("B7").Value = 28        # Winning Trades
$summary.Range("B9").Value = 33.73     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98999999999999   # Capital
$status.Range("D4").Value = 83                  # Trades
$status.Range("E4").Value = -0.02               # P&L $
$status.Range("F4").Value = -0.01               # P&L %
$status.Range("G4").Value = 33.73               # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the new trade #83 row to a trades-log sheet at row 84.
# Date/Time-like text must stay text (matching the log's plain-text
# columns), so a quote-prefix is used to stop Excel's auto date/time
# parsing from converting them into date serials.
# ---------------------------------------------------------------------------
function Add-Trade83Row($ws) {
    $ws.Cells.Item(84, 1).Value  = 83
    $ws.Cells.Item(84, 2).Value  = "'2026-02-17"
    $ws.Cells.Item(84, 2).Style  = "Normal"
    $ws.Cells.Item(84, 3).Value  = "'15:53:00"
    $ws.Cells.Item(84, 3).Style  = "Normal"
    $ws.Cells.Item(84, 4).Value  = "MarketMaking"
    $ws.Cells.Item(84, 5).Value  = "DOWN"
    $ws.Cells.Item(84, 6).Value  = 0.17
    $ws.Cells.Item(84, 7).Value  = 0.2
    $ws.Cells.Item(84, 8).Value  = "CLOSED"
    $ws.Cells.Item(84, 9).Value  = 17.6471
    $ws.Cells.Item(84, 10).Value = 0.03
    $ws.Cells.Item(84, 11).Value = 99.98999999999999
    $ws.Cells.Item(84, 12).Value = 0
    $ws.Cells.Item(84, 13).Value = 0
    $ws.Cells.Item(84, 14).Value = 0.6
    $ws.Cells.Item(84, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(84, 16).Value = "early_exit"
    $ws.Cells.Item(84, 17).Value = 0.13
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade83Row $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade83Row $marketMaking

Write-Host "Applied trade #83 update to Summary, Strategy Status, All Trades, MarketMaking"
